# Update "想去人数" (number of people interested) counts on the
# "展览" and "全部类型" worksheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 1330
    $ws.Range("F3").Value = 1838
    $ws.Range("F4").Value = 144
    $ws.Range("F7").Value = 161
}
